$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve leading zeros / text formatting for the Codice Fiscale and count columns,
# matching the original inline-string cell types in the source file.
$ws.Columns.Item(2).NumberFormat = "@"
$ws.Columns.Item(3).NumberFormat = "@"

$ws.Cells.Item(3, 1).Value = "Unione Italiana delle Camere di Commercio Industria, Artigianato e Agricoltura"
$ws.Cells.Item(3, 2).Value = "01484460587"
$ws.Cells.Item(3, 3).Value = "1080"
$ws.Cells.Item(6, 1).Value = "Regione Lombardia"
$ws.Cells.Item(6, 2).Value = "80050050154"
$ws.Cells.Item(6, 3).Value = "504"
$ws.Cells.Item(7, 1).Value = "Maggioli SPA"
$ws.Cells.Item(7, 2).Value = "06188330150"
$ws.Cells.Item(7, 3).Value = "440"
$ws.Cells.Item(9, 1).Value = "CREDEMTEL SpA"
$ws.Cells.Item(9, 2).Value = "01378570350"
$ws.Cells.Item(9, 3).Value = "416"
$ws.Cells.Item(11, 1).Value = "Banca Popolare di Sondrio, Società Cooperativa per Azioni"
$ws.Cells.Item(11, 2).Value = "00053810149"
$ws.Cells.Item(11, 3).Value = "351"
$ws.Cells.Item(14, 1).Value = "APKAPPA S.R.L."
$ws.Cells.Item(14, 2).Value = "08543640158"
$ws.Cells.Item(14, 3).Value = "241"
$ws.Cells.Item(16, 1).Value = "Progetti e Soluzioni SPA"
$ws.Cells.Item(16, 2).Value = "06423240727"
$ws.Cells.Item(16, 3).Value = "198"
$ws.Cells.Item(22, 1).Value = "Regione Puglia"
$ws.Cells.Item(22, 2).Value = "80017210727"
$ws.Cells.Item(22, 3).Value = "125"
$ws.Cells.Item(26, 1).Value = "Siscom SPA"
$ws.Cells.Item(26, 2).Value = "01778000040"
$ws.Cells.Item(26, 3).Value = "98"
$ws.Cells.Item(27, 1).Value = "Regione Autonoma Friuli-Venezia Giulia"
$ws.Cells.Item(27, 2).Value = "80014930327"
$ws.Cells.Item(27, 3).Value = "98"
$ws.Cells.Item(29, 1).Value = "PMPay s.r.l."
$ws.Cells.Item(29, 2).Value = "08747230962"
$ws.Cells.Item(29, 3).Value = "77"
$ws.Cells.Item(31, 1).Value = "Intesa Sanpaolo SPA"
$ws.Cells.Item(31, 2).Value = "00799960158"
$ws.Cells.Item(31, 3).Value = "68"
$ws.Cells.Item(34, 1).Value = "Regione Umbria"
$ws.Cells.Item(34, 2).Value = "80000130544"
$ws.Cells.Item(34, 3).Value = "55"
$ws.Cells.Item(42, 1).Value = "Regione Liguria"
$ws.Cells.Item(42, 2).Value = "00849050109"
$ws.Cells.Item(42, 3).Value = "34"
$ws.Cells.Item(44, 1).Value = "Numera Sistemi e Informatica SpA"
$ws.Cells.Item(44, 2).Value = "01265230902"
$ws.Cells.Item(44, 3).Value = "30"
$ws.Cells.Item(51, 1).Value = "Si.Form Consulting srl"
$ws.Cells.Item(51, 2).Value = "03943960827"
$ws.Cells.Item(51, 3).Value = "15"
$ws.Cells.Item(52, 1).Value = "Servizi Locali SpA"
$ws.Cells.Item(52, 2).Value = "03170580751"
$ws.Cells.Item(52, 3).Value = "15"
$ws.Cells.Item(55, 1).Value = "UBI Banca"
$ws.Cells.Item(55, 2).Value = "03053920165"
$ws.Cells.Item(55, 3).Value = "8"
$ws.Cells.Item(56, 1).Value = "Comune di Catania"
$ws.Cells.Item(56, 2).Value = "00137020871"
$ws.Cells.Item(56, 3).Value = "8"
$ws.Cells.Item(60, 1).Value = "ARGO SOFTWARE SRL"
$ws.Cells.Item(60, 2).Value = "00838520880"
$ws.Cells.Item(60, 3).Value = "4"
$ws.Cells.Item(61, 1).Value = "Phoenix IT Solutions S.r.L"
$ws.Cells.Item(61, 2).Value = "07623321218"
$ws.Cells.Item(61, 3).Value = "4"
$ws.Cells.Item(66, 1).Value = "KOINE' SRL"
$ws.Cells.Item(66, 2).Value = "01934790971"
$ws.Cells.Item(66, 3).Value = "2"
$ws.Cells.Item(67, 1).Value = "ICCREA Banca SpA"
$ws.Cells.Item(67, 2).Value = "04774801007"
$ws.Cells.Item(67, 3).Value = "2"
$ws.Cells.Item(68, 1).Value = "I.C.A. - Imposte Comunali Affini – s.r.l."
$ws.Cells.Item(68, 2).Value = "02478610583"
$ws.Cells.Item(68, 3).Value = "1"
$ws.Cells.Item(69, 1).Value = "Agenzia Italiana del Farmaco - AIFA"
$ws.Cells.Item(69, 2).Value = "97345810580"
$ws.Cells.Item(69, 3).Value = "1"
$ws.Cells.Item(70, 1).Value = "Banco BPM Società per Azioni"
$ws.Cells.Item(70, 2).Value = "09722490969"
$ws.Cells.Item(70, 3).Value = "1"
$ws.Cells.Item(71, 1).Value = "Ministero dello Sviluppo Economico"
$ws.Cells.Item(71, 2).Value = "80230390587"
$ws.Cells.Item(71, 3).Value = "1"
$ws.Cells.Item(72, 1).Value = "Noviservice srl"
$ws.Cells.Item(72, 2).Value = "02789990922"
$ws.Cells.Item(72, 3).Value = "1"
$ws.Cells.Item(73, 1).Value = "Softline srl"
$ws.Cells.Item(73, 2).Value = "12299030150"
$ws.Cells.Item(73, 3).Value = "1"
$ws.Cells.Item(74, 1).Value = "MegASP S.r.l."
$ws.Cells.Item(74, 2).Value = "09898030151"
$ws.Cells.Item(74, 3).Value = "1"
$ws.Cells.Item(75, 1).Value = "Banca Nazionale del Lavoro S.p.A."
$ws.Cells.Item(75, 2).Value = "09339391006"
$ws.Cells.Item(75, 3).Value = "1"
$ws.Cells.Item(76, 1).Value = "Engineering Ingegneria Informatica SpA"
$ws.Cells.Item(76, 2).Value = "00967720285"
$ws.Cells.Item(76, 3).Value = "1"
$ws.Cells.Item(77, 1).Value = "BANCA MONTE DEI PASCHI DI SIENA"
$ws.Cells.Item(77, 2).Value = "00884060526"
$ws.Cells.Item(77, 3).Value = "1"
$ws.Cells.Item(79, 1).Value = "Società Almaviva S.p.A."
$ws.Cells.Item(79, 2).Value = "08450891000"
$ws.Cells.Item(79, 3).Value = "1"
